{"js": "// Replace each two-digit multiplication expression in the document's single\n// table with its new value, per the commit's regenerated numbers.\n// Each \"find\" string below is unique within the document, so a plain\n// search-and-replace (matchCase, not whole-word since \"\u00d7\" and \"=\" already\n// bound the token) safely targets exactly one cell each.\nconst replacements = [\n  [\"16\u00d759=\", \"35\u00d746=\"],\n  [\"79\u00d756=\", \"92\u00d769=\"],\n  [\"68\u00d760=\", \"28\u00d715=\"],\n  [\"71\u00d735=\", \"19\u00d769=\"],\n  [\"17\u00d766=\", \"51\u00d776=\"],\n  [\"42\u00d788=\", \"47\u00d788=\"],\n  [\"60\u00d750=\", \"39\u00d740=\"],\n  [\"98\u00d767=\", \"17\u00d787=\"],\n  [\"79\u00d761=\", \"39\u00d746=\"],\n  [\"68\u00d777=\", \"58\u00d786=\"],\n  [\"97\u00d741=\", \"79\u00d715=\"],\n  [\"70\u00d795=\", \"94\u00d797=\"],\n  [\"82\u00d775=\", \"42\u00d721=\"],\n  [\"23\u00d735=\", \"72\u00d728=\"],\n  [\"35\u00d777=\", \"51\u00d769=\"],\n  [\"94\u00d725=\", \"55\u00d718=\"],\n  [\"44\u00d730=\", \"52\u00d744=\"],\n  [\"79\u00d718=\", \"35\u00d775=\"],\n  [\"81\u00d788=\", \"41\u00d791=\"],\n  [\"32\u00d792=\", \"57\u00d773=\"],\n  [\"54\u00d712=\", \"48\u00d767=\"],\n  [\"53\u00d766=\", \"38\u00d760=\"],\n  [\"31\u00d729=\", \"95\u00d783=\"],\n  [\"74\u00d789=\", \"80\u00d724=\"],\n  [\"73\u00d792=\", \"77\u00d772=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication expression in the worksheet's table\n# to the newly generated value. Every \"find\" string is unique in the\n# document, so Find/Replace with ReplaceAll (one hit each) safely retargets\n# exactly the intended cell without touching anything else.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"16\u00d759=\", \"35\u00d746=\"),\n    @(\"79\u00d756=\", \"92\u00d769=\"),\n    @(\"68\u00d760=\", \"28\u00d715=\"),\n    @(\"71\u00d735=\", \"19\u00d769=\"),\n    @(\"17\u00d766=\", \"51\u00d776=\"),\n    @(\"42\u00d788=\", \"47\u00d788=\"),\n    @(\"60\u00d750=\", \"39\u00d740=\"),\n    @(\"98\u00d767=\", \"17\u00d787=\"),\n    @(\"79\u00d761=\", \"39\u00d746=\"),\n    @(\"68\u00d777=\", \"58\u00d786=\"),\n    @(\"97\u00d741=\", \"79\u00d715=\"),\n    @(\"70\u00d795=\", \"94\u00d797=\"),\n    @(\"82\u00d775=\", \"42\u00d721=\"),\n    @(\"23\u00d735=\", \"72\u00d728=\"),\n    @(\"35\u00d777=\", \"51\u00d769=\"),\n    @(\"94\u00d725=\", \"55\u00d718=\"),\n    @(\"44\u00d730=\", \"52\u00d744=\"),\n    @(\"79\u00d718=\", \"35\u00d775=\"),\n    @(\"81\u00d788=\", \"41\u00d791=\"),\n    @(\"32\u00d792=\", \"57\u00d773=\"),\n    @(\"54\u00d712=\", \"48\u00d767=\"),\n    @(\"53\u00d766=\", \"38\u00d760=\"),\n    @(\"31\u00d729=\", \"95\u00d783=\"),\n    @(\"74\u00d789=\", \"80\u00d724=\"),\n    @(\"73\u00d792=\", \"77\u00d772=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1  # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute([ref]$find, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replace, [ref]2) | Out-Null\n}\n"}
